# Applies the edit described by the diff:
#  - Sheet "1碑影迷踪": cell B2 loses its trailing full-width period.
#  - The active/selected worksheet moves from "2消失的龙" to "1碑影迷踪".
#  - On "1碑影迷踪" the selection moves from B2 to B5.
#  - On "2消失的龙" the sheet is no longer the selected/active tab.

$wb = $excel.ActiveWorkbook

$wsInfo = $wb.Worksheets.Item("1碑影迷踪")
$wsLink = $wb.Worksheets.Item("2消失的龙")

# Strip the trailing "。" from the prompt text in B2 (same text, just no
# trailing full stop).
$wsInfo.Range("B2").Value = "日记中的古诗提到两处古迹的名字，有一处现已无存，请找到仍然存在的古迹名称"

# Make "1碑影迷踪" the active sheet/tab and move its selection to B5.
$wsInfo.Activate() | Out-Null
$wsInfo.Range("B5").Select() | Out-Null

# Keep "2消失的龙" selection at C8 (unchanged) - just no longer the active tab.
$wsLink.Range("C8").Select() | Out-Null

# Re-activate "1碑影迷踪" so it remains the active sheet on save.
$wsInfo.Activate() | Out-Null
$wsInfo.Range("B5").Select() | Out-Null
